{"js": "// Update the worksheet date header and all 25 three-digit x one-digit\n// multiplication equations to the next day's values.\n// Every \"before\" string is unique in the document, so an exact\n// (matchCase) search-and-replace per pair is unambiguous.\nconst replacements = [\n  [\"2024-09-03 Tuesday\", \"2024-09-04 Wednesday\"],\n  [\"545\u00d73=1635\", \"285\u00d75=1425\"],\n  [\"507\u00d79=4563\", \"864\u00d79=7776\"],\n  [\"912\u00d78=7296\", \"615\u00d78=4920\"],\n  [\"332\u00d79=2988\", \"211\u00d74=844\"],\n  [\"259\u00d78=2072\", \"855\u00d75=4275\"],\n  [\"142\u00d72=284\", \"462\u00d74=1848\"],\n  [\"641\u00d74=2564\", \"415\u00d76=2490\"],\n  [\"669\u00d76=4014\", \"983\u00d75=4915\"],\n  [\"348\u00d72=696\", \"213\u00d79=1917\"],\n  [\"158\u00d73=474\", \"428\u00d78=3424\"],\n  [\"318\u00d78=2544\", \"567\u00d78=4536\"],\n  [\"540\u00d74=2160\", \"718\u00d76=4308\"],\n  [\"608\u00d78=4864\", \"194\u00d76=1164\"],\n  [\"713\u00d76=4278\", \"401\u00d78=3208\"],\n  [\"203\u00d79=1827\", \"780\u00d76=4680\"],\n  [\"624\u00d76=3744\", \"841\u00d77=5887\"],\n  [\"859\u00d77=6013\", \"554\u00d76=3324\"],\n  [\"718\u00d73=2154\", \"998\u00d77=6986\"],\n  [\"709\u00d76=4254\", \"717\u00d76=4302\"],\n  [\"204\u00d76=1224\", \"221\u00d76=1326\"],\n  [\"533\u00d74=2132\", \"135\u00d76=810\"],\n  [\"152\u00d74=608\", \"823\u00d75=4115\"],\n  [\"210\u00d73=630\", \"820\u00d76=4920\"],\n  [\"947\u00d75=4735\", \"442\u00d72=884\"],\n  [\"707\u00d72=1414\", \"775\u00d75=3875\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date header and all 25 three-digit x one-digit\n# multiplication equations to the next day's values.\n# Every \"before\" string is unique in the document, so an exact\n# (MatchCase) Find/Replace per pair is unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-09-03 Tuesday\", \"2024-09-04 Wednesday\"),\n    @(\"545\u00d73=1635\", \"285\u00d75=1425\"),\n    @(\"507\u00d79=4563\", \"864\u00d79=7776\"),\n    @(\"912\u00d78=7296\", \"615\u00d78=4920\"),\n    @(\"332\u00d79=2988\", \"211\u00d74=844\"),\n    @(\"259\u00d78=2072\", \"855\u00d75=4275\"),\n    @(\"142\u00d72=284\", \"462\u00d74=1848\"),\n    @(\"641\u00d74=2564\", \"415\u00d76=2490\"),\n    @(\"669\u00d76=4014\", \"983\u00d75=4915\"),\n    @(\"348\u00d72=696\", \"213\u00d79=1917\"),\n    @(\"158\u00d73=474\", \"428\u00d78=3424\"),\n    @(\"318\u00d78=2544\", \"567\u00d78=4536\"),\n    @(\"540\u00d74=2160\", \"718\u00d76=4308\"),\n    @(\"608\u00d78=4864\", \"194\u00d76=1164\"),\n    @(\"713\u00d76=4278\", \"401\u00d78=3208\"),\n    @(\"203\u00d79=1827\", \"780\u00d76=4680\"),\n    @(\"624\u00d76=3744\", \"841\u00d77=5887\"),\n    @(\"859\u00d77=6013\", \"554\u00d76=3324\"),\n    @(\"718\u00d73=2154\", \"998\u00d77=6986\"),\n    @(\"709\u00d76=4254\", \"717\u00d76=4302\"),\n    @(\"204\u00d76=1224\", \"221\u00d76=1326\"),\n    @(\"533\u00d74=2132\", \"135\u00d76=810\"),\n    @(\"152\u00d74=608\", \"823\u00d75=4115\"),\n    @(\"210\u00d73=630\", \"820\u00d76=4920\"),\n    @(\"947\u00d75=4735\", \"442\u00d72=884\"),\n    @(\"707\u00d72=1414\", \"775\u00d75=3875\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n\nWrite-Output \"done\"\n"}
